$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated
# from 45186 to 45188 for every data row (rows 2 through 540).
$ws.Range("C2:C540").Value = 45188
